$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Re-style the existing "Attachment" metadata row (row 3): D3:F3 pick up
#    the plain default font (previously unstyled / style 0) and G3 (the
#    modification-date cell) switches from the date-formatted style to a
#    bold-black General style.
# ---------------------------------------------------------------------------
$ws.Range("D3:F3").Font.Name = "Arial"
$ws.Range("D3:F3").Font.Name = "Calibri"

$ws.Range("G3").NumberFormat = "General"
$ws.Range("G3").Font.Bold = $true
$ws.Range("G3").Font.Size = 11
$ws.Range("G3").Font.Name = "Calibri"
$ws.Range("G3").Font.Color = 0

# ---------------------------------------------------------------------------
# 2) Add the two new property-flag columns: "Multivalued" and "Unique".
#    Headers (row 4) use the bold black Calibri header font; the data rows
#    (5-7) hold literal text "FALSE" formatted with a TRUE/FALSE custom
#    number format and left alignment. Row 8 gets a single styled, empty
#    anchor cell in column K (mirrors the sheet's extended used range).
# ---------------------------------------------------------------------------
$ws.Range("K4").Value = "Multivalued"
$ws.Range("L4").Value = "Unique"
$ws.Range("K4:L4").Font.Bold = $true
$ws.Range("K4:L4").Font.Size = 11
$ws.Range("K4:L4").Font.Name = "Calibri"
$ws.Range("K4:L4").Font.Color = 0

$ws.Range("K5").Value = "'FALSE"
$ws.Range("L5").Value = "'FALSE"
$ws.Range("K6").Value = "'FALSE"
$ws.Range("L6").Value = "'FALSE"
$ws.Range("K7").Value = "'FALSE"
$ws.Range("L7").Value = "'FALSE"

$ws.Range("K5:L7").NumberFormat = """TRUE"";""TRUE"";""FALSE"""
$ws.Range("K5:L7").HorizontalAlignment = -4131

$ws.Range("K8").NumberFormat = """TRUE"";""TRUE"";""FALSE"""
$ws.Range("K8").HorizontalAlignment = -4131
$ws.Rows(8).RowHeight = 15

# ---------------------------------------------------------------------------
# 3) Selection / active cell mirrors what the author left selected after the
#    edit.
# ---------------------------------------------------------------------------
$ws.Range("L4:L7").Select()
$ws.Application.ActiveCell = $ws.Range("L4")
